# Actualización automática hashcode
# Updates the MD5 hash values in column B for a set of rows identified by
# their unique "codigo" key in column A (see mapping comments below).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 05-050301A
$ws.Range("B11").Value = "2b90794bbf410b78bd68be02a9afdc28"
# 05-050302A
$ws.Range("B29").Value = "4ec0f949f570e44f49b8f9d5ec6e1c20"
# 05-050301TP
$ws.Range("B121").Value = "663ee48c075b606ef9eb43f19f08fa8d"
# 05-050301TC
$ws.Range("B123").Value = "b870e4d23cc1caa8e658f3ef81ba8eb8"
# 05-050312TP
$ws.Range("B133").Value = "97aeef0afb4f48bac25c35d5a8352971"
# 05-050203TP
$ws.Range("B159").Value = "1fbda8d6ff8792b063dccee95965b508"
# 05-050314TC
$ws.Range("B198").Value = "d8de88e2e28fb88894f2abd73abd6529"
# 05-050003TC
$ws.Range("B246").Value = "523ee4c716a1c651a24a160e0173d9fc"
# 05-050003TP
$ws.Range("B279").Value = "0f715dab3fecfd2e04f6b6803d5bac9f"
# 05-050001TC
$ws.Range("B414").Value = "40cef0dea6b96d7d65a86d830dab5a11"
# 05-0709-070905BTC
$ws.Range("B423").Value = "930e9bd628ccd09c643cd2b4a4b8cfad"
# 05-050001TP
$ws.Range("B451").Value = "7284fed6f381b854c6cf32dc28a30074"
# 05-050312A
$ws.Range("B525").Value = "7adb1e39b82cd9c8011a353bdbaab39f"
# 05-050004A
$ws.Range("B578").Value = "bf877eeb2e688c5f25b5113f5e4ad1f7"
# 05-050204TP
$ws.Range("B628").Value = "31ab3308ae96077bf0b4424dc57cdc95"
# 05-050204TC
$ws.Range("B639").Value = "315916cfa35efed5711d9559b937f838"
# 05-050302TP
$ws.Range("B641").Value = "c85a0212c77cacd97ca482f471a84fe4"
# 05-050313TP
$ws.Range("B661").Value = "0068163e6ab0852da20745f8fc355361"
# 05-050004TC
$ws.Range("B768").Value = "979626509ac5420633d6882d02ffb4e7"
# 05-050004TP
$ws.Range("B773").Value = "d2251d4d8c886493d6a8b47207eede45"
# 05-050002TC
$ws.Range("B798").Value = "db9d40ce772e375861da826f281c42dc"
# 05-050003A
$ws.Range("B881").Value = "d37e25fb56cb04197ac800a229c553fd"
# 05-050001A
$ws.Range("B917").Value = "6bfb93d48c5d0590e1b5a7594ffdda98"
# 05-050002TP
$ws.Range("B941").Value = "14d359d1baf1e0f5985c646732c37f21"
